# The deck ships two DrawingML themes:
#   theme1.xml -> applied theme ("Integral" / "Red Violet" colours)
#   theme2.xml -> the stock "Office Theme" colours (only wired to the
#                 notes master in this file, not reachable as a Design)
#
# The authored change swaps the two themes' colour schemes: the slide
# master's applied theme becomes the standard "Office Theme" palette.
# We recolour the live design's ThemeColorScheme to match that palette,
# one ThemeColorSchemeIndex slot at a time (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), using the classic
# RGB = blue*65536 + green*256 + red packing that PowerPoint's object
# model expects for a ThemeColor.RGB value.

$p  = $ppt.ActivePresentation
$d  = $p.Designs.Item(1)
$sm = $d.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

$cs.Colors(1).RGB  = 0x00*65536 + 0x00*256 + 0x00   # dk1      -> 000000
$cs.Colors(2).RGB  = 0xFF*65536 + 0xFF*256 + 0xFF   # lt1      -> FFFFFF
$cs.Colors(3).RGB  = 0x6A*65536 + 0x54*256 + 0x44   # dk2      -> 44546A
$cs.Colors(4).RGB  = 0xE6*65536 + 0xE6*256 + 0xE7   # lt2      -> E7E6E6
$cs.Colors(5).RGB  = 0xD5*65536 + 0x9B*256 + 0x5B   # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = 0x31*65536 + 0x7D*256 + 0xED   # accent2  -> ED7D31
$cs.Colors(7).RGB  = 0xA5*65536 + 0xA5*256 + 0xA5   # accent3  -> A5A5A5
$cs.Colors(8).RGB  = 0x00*65536 + 0xC0*256 + 0xFF   # accent4  -> FFC000
$cs.Colors(9).RGB  = 0xC4*65536 + 0x72*256 + 0x44   # accent5  -> 4472C4
$cs.Colors(10).RGB = 0x47*65536 + 0xAD*256 + 0x70   # accent6  -> 70AD47
$cs.Colors(11).RGB = 0xC1*65536 + 0x63*256 + 0x05   # hlink    -> 0563C1
$cs.Colors(12).RGB = 0x72*65536 + 0x4F*256 + 0x95   # folHlink -> 954F72
